# Add a "Type" column (E) to the WiscSIMSColumnDictionary sheet, designating
# each data row as either "Text" or "Numeric".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("E1").Value = "Type"

# Row -> Type mapping for data rows 2..39
$types = @{
    2  = "Text"
    3  = "Text"
    4  = "Numeric"
    5  = "Numeric"
    6  = "Numeric"
    7  = "Numeric"
    8  = "Numeric"
    9  = "Numeric"
    10 = "Numeric"
    11 = "Numeric"
    12 = "Numeric"
    13 = "Text"
    14 = "Text"
    15 = "Numeric"
    16 = "Numeric"
    17 = "Numeric"
    18 = "Numeric"
    19 = "Numeric"
    20 = "Numeric"
    21 = "Numeric"
    22 = "Numeric"
    23 = "Numeric"
    24 = "Numeric"
    25 = "Numeric"
    26 = "Text"
    27 = "Numeric"
    28 = "Text"
    29 = "Text"
    30 = "Text"
    31 = "Text"
    32 = "Numeric"
    33 = "Numeric"
    34 = "Numeric"
    35 = "Numeric"
    36 = "Text"
    37 = "Text"
    38 = "Text"
    39 = "Numeric"
}

foreach ($row in 2..39) {
    $ws.Range("E$row").Value = $types[$row]
}

# Rows 15-39 carry an explicit black font colour on the new cells.
$ws.Range("E15:E39").Font.Color = 0

# Update the view: move the selection to E39 (matches the saved workbook state).
$ws.Range("E39").Select()
